# Update column F (dSF) values for the specified rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -8
    3  = -8
    4  = -4
    5  = -10
    14 = -2
    16 = -5
    18 = -5
    19 = -1
    39 = 2
    41 = -5
    46 = -4
    48 = -1
    50 = -2
    51 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
